# Updates cryptos list values (Price and Volume(1h) columns) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Prefix with an apostrophe so Excel always stores the value as literal text
    # (never auto-converts numeric-looking strings like "532.17" to a Number),
    # then restore the 'Normal' style so no quote-prefix formatting is retained.
    $range.Value = "'" + $value
    $range.Style = 'Normal'
}

Set-TextValue $ws.Range('D2') '57.765.28'
Set-TextValue $ws.Range('E2') '  -3.42%  '
Set-TextValue $ws.Range('D3') '2.287.56'
Set-TextValue $ws.Range('E3') '  -3.61%  '
Set-TextValue $ws.Range('E4') '  -0.08%  '
Set-TextValue $ws.Range('D5') '532.17'
Set-TextValue $ws.Range('E5') '  -4.25%  '
Set-TextValue $ws.Range('D6') '130.35'
Set-TextValue $ws.Range('E6') '  -2.36%  '
Set-TextValue $ws.Range('E7') '  -0.03%  '
Set-TextValue $ws.Range('D8') '0.579'
Set-TextValue $ws.Range('E8') '  -1.28%  '
Set-TextValue $ws.Range('D9') '2.288.48'
Set-TextValue $ws.Range('E9') '  -3.42%  '
Set-TextValue $ws.Range('D10') '0.0993'
Set-TextValue $ws.Range('E10') '  -5.54%  '
Set-TextValue $ws.Range('D11') '5.42'
Set-TextValue $ws.Range('E11') '  -4.18%  '
Set-TextValue $ws.Range('E12') '  -0.56%  '
Set-TextValue $ws.Range('D13') '0.329'
Set-TextValue $ws.Range('E13') '  -3.86%  '
Set-TextValue $ws.Range('D14') '23.40'
Set-TextValue $ws.Range('E14') '  -4.07%  '
Set-TextValue $ws.Range('D15') '2.698.84'
Set-TextValue $ws.Range('E15') '  -3.65%  '
Set-TextValue $ws.Range('D16') '57.753.18'
Set-TextValue $ws.Range('E16') '  -3.36%  '
Set-TextValue $ws.Range('E17') '  -4.41%  '
Set-TextValue $ws.Range('D18') '2.290.37'
Set-TextValue $ws.Range('E18') '  -3.48%  '
Set-TextValue $ws.Range('E19') '  -5.47%  '
Set-TextValue $ws.Range('E20') '  -5.78%  '
Set-TextValue $ws.Range('D21') '312.62'
Set-TextValue $ws.Range('E21') '  -2.68%  '
Set-TextValue $ws.Range('D22') '6.34'
Set-TextValue $ws.Range('E22') '  -4.43%  '
Set-TextValue $ws.Range('E23') '  -0.01%  '
Set-TextValue $ws.Range('D24') '62.37'
Set-TextValue $ws.Range('E24') '  -2.73%  '
Set-TextValue $ws.Range('E25') '  -4.86%  '
Set-TextValue $ws.Range('D26') '0.999'
Set-TextValue $ws.Range('E26') '  -0.15%  '
Set-TextValue $ws.Range('E27') '  -4.48%  '
Set-TextValue $ws.Range('D28') '1.27'
Set-TextValue $ws.Range('E28') '  -6.64%  '
Set-TextValue $ws.Range('D29') '170.64'
Set-TextValue $ws.Range('E29') '  +0.33%  '
Set-TextValue $ws.Range('E30') '  -5.65%  '
Set-TextValue $ws.Range('D31') '0.0₃0714'
Set-TextValue $ws.Range('E31') '  -5.70%  '
Set-TextValue $ws.Range('D32') '5.74'
Set-TextValue $ws.Range('E32') '  -5.22%  '
Set-TextValue $ws.Range('E33') '  -6.04%  '
Set-TextValue $ws.Range('D34') '0.379'
Set-TextValue $ws.Range('E34') '  -4.89%  '
Set-TextValue $ws.Range('E35') '  +0.03%  '
Set-TextValue $ws.Range('D36') '17.68'
Set-TextValue $ws.Range('E36') '  -2.59%  '
Set-TextValue $ws.Range('E37') '  -0.02%  '
Set-TextValue $ws.Range('E38') '  -7.07%  '
Set-TextValue $ws.Range('D39') '3.88'
Set-TextValue $ws.Range('E39') '  -6.14%  '
Set-TextValue $ws.Range('D40') '38.09'
Set-TextValue $ws.Range('E40') '  -1.30%  '
Set-TextValue $ws.Range('E41') '  -6.24%  '
Set-TextValue $ws.Range('D42') '140.89'
Set-TextValue $ws.Range('E42') '  -2.76%  '
Set-TextValue $ws.Range('D43') '287.59'
Set-TextValue $ws.Range('E43') '  -9.54%  '
Set-TextValue $ws.Range('D44') '3.39'
Set-TextValue $ws.Range('E44') '  -3.86%  '
Set-TextValue $ws.Range('D45') '0.0946'
Set-TextValue $ws.Range('E45') '  -2.37%  '
Set-TextValue $ws.Range('E46') '  -2.66%  '
Set-TextValue $ws.Range('E47') '  -3.04%  '
Set-TextValue $ws.Range('D48') '18.03'
Set-TextValue $ws.Range('E48') '  -8.11%  '
Set-TextValue $ws.Range('D49') '0.0209'
Set-TextValue $ws.Range('E49') '  -3.58%  '
Set-TextValue $ws.Range('E50') '  -1.09%  '
Set-TextValue $ws.Range('D51') '0.0₆0201'
Set-TextValue $ws.Range('E51') '  +84.89%  '
